$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1058.0588
$ws.Range("I15").Value = 1058.0588
$ws.Range("K15").Value = 3174.1764
$ws.Range("M15").Value = -3005.1764

$ws.Range("H33").Value = 5883699.5
$ws.Range("I33").Value = 8333491
$ws.Range("J33").Value = 4200
$ws.Range("K33").Value = 8333491
$ws.Range("L33").Value = 4200
$ws.Range("M33").Value = -8333262
$ws.Range("N33").Value = -4658

$ws.Range("H76").Value = 4000.5454
$ws.Range("I76").Value = 3778
$ws.Range("K76").Value = 3778
$ws.Range("M76").Value = -3463

$ws.Range("H79").Value = 4000.5454
$ws.Range("I79").Value = 3778
$ws.Range("K79").Value = 3778
$ws.Range("M79").Value = -2686

$ws.Range("H96").Value = 269.9
$ws.Range("J96").Value = 130
$ws.Range("L96").Value = 390
$ws.Range("N96").Value = -3136

$ws.Range("H99").Value = 2681.8948
$ws.Range("I99").Value = 919.7143
$ws.Range("J99").Value = 3709.8333
$ws.Range("K99").Value = 2759.1429
$ws.Range("L99").Value = 11129.4999
$ws.Range("M99").Value = -1261.1429
$ws.Range("N99").Value = -14125.4999

$ws.Range("H137").Value = 1776762.1
$ws.Range("I137").Value = 25051.191
$ws.Range("K137").Value = 75153.573
$ws.Range("M137").Value = -72603.573

$ws.Range("H138").Value = 4458.674
$ws.Range("I138").Value = 2377.9285
$ws.Range("J138").Value = 4832.141
$ws.Range("K138").Value = 7133.7855
$ws.Range("L138").Value = 14496.423
$ws.Range("M138").Value = -1993.7855
$ws.Range("N138").Value = -24776.423

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 13097
$ws.Range("I22").Value = 155.4
$ws.Range("J22").Value = 34666.332
$ws.Range("K22").Value = 155.4
$ws.Range("L22").Value = 34666.332
$ws.Range("M22").Value = 143.6
$ws.Range("N22").Value = -35264.332

$ws.Range("H32").Value = 2239.1973
$ws.Range("I32").Value = 1552.3881
$ws.Range("J32").Value = 13743.25
$ws.Range("K32").Value = 1552.3881
$ws.Range("L32").Value = 13743.25
$ws.Range("M32").Value = -1265.3881
$ws.Range("N32").Value = -14317.25

$ws.Range("H61").Value = 2743.5557
$ws.Range("I61").Value = 2447.4614
$ws.Range("J61").Value = 3513.4
$ws.Range("K61").Value = 2447.4614
$ws.Range("L61").Value = 3513.4
$ws.Range("M61").Value = -2235.4614
$ws.Range("N61").Value = -3937.4

$ws.Range("H74").Value = 29445900
$ws.Range("I74").Value = 37496.535
$ws.Range("J74").Value = 166685120
$ws.Range("K74").Value = 37496.535
$ws.Range("L74").Value = 166685120
$ws.Range("M74").Value = -36622.535
$ws.Range("N74").Value = -166686868

$ws.Range("H77").Value = 29445900
$ws.Range("I77").Value = 37496.535
$ws.Range("J77").Value = 166685120
$ws.Range("K77").Value = 187482.675
$ws.Range("L77").Value = 833425600
$ws.Range("M77").Value = -183114.675
$ws.Range("N77").Value = -833434336

$ws.Range("H132").Value = 3610.2917
$ws.Range("I132").Value = 3341.3076
$ws.Range("K132").Value = 10023.9228
$ws.Range("M132").Value = -7493.9228

$ws.Range("H136").Value = 2743.5557
$ws.Range("I136").Value = 2447.4614
$ws.Range("J136").Value = 3513.4
$ws.Range("K136").Value = 7342.3842
$ws.Range("L136").Value = 10540.2
$ws.Range("M136").Value = -4792.3842
$ws.Range("N136").Value = -15640.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2053.8928
$ws.Range("I134").Value = 1824.36
$ws.Range("K134").Value = 5473.08
$ws.Range("M134").Value = -2938.08

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23476052
$ws.Range("I31").Value = 3726750.2
$ws.Range("J31").Value = 71438650
$ws.Range("K31").Value = 3726750.2
$ws.Range("L31").Value = 71438650
$ws.Range("M31").Value = -3726455.2
$ws.Range("N31").Value = -71439240

$ws.Range("H34").Value = 23476052
$ws.Range("I34").Value = 3726750.2
$ws.Range("J34").Value = 71438650
$ws.Range("K34").Value = 3726750.2
$ws.Range("L34").Value = 71438650
$ws.Range("M34").Value = -3726548.2
$ws.Range("N34").Value = -71439054

$ws.Range("H122").Value = 1990
$ws.Range("J122").Value = 1990
$ws.Range("L122").Value = 5970
$ws.Range("N122").Value = -10870

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3539.8333
$ws.Range("J5").Value = 5958.077
$ws.Range("L5").Value = 17874.231
$ws.Range("N5").Value = -18098.231

$ws.Range("H86").Value = 301.72726
$ws.Range("I86").Value = 60
$ws.Range("K86").Value = 180
$ws.Range("M86").Value = 1006

$ws.Range("H89").Value = 301.72726
$ws.Range("I89").Value = 60
$ws.Range("K89").Value = 540
$ws.Range("M89").Value = 5388

$ws.Range("H113").Value = 377.37036
$ws.Range("I113").Value = 356.2353
$ws.Range("J113").Value = 413.3
$ws.Range("K113").Value = 1068.7059
$ws.Range("L113").Value = 1239.9
$ws.Range("M113").Value = 1101.2941
$ws.Range("N113").Value = -5579.9

$ws.Range("H118").Value = 2559.2
$ws.Range("I118").Value = 2559.2
$ws.Range("K118").Value = 7677.599999999999
$ws.Range("M118").Value = -6434.599999999999

$ws.Range("H122").Value = 1052.0769
$ws.Range("I122").Value = 519.1429000000001
$ws.Range("J122").Value = 1673.8334
$ws.Range("K122").Value = 4672.2861
$ws.Range("L122").Value = 15064.5006
$ws.Range("M122").Value = -2222.2861
$ws.Range("N122").Value = -19964.5006

$ws.Range("H124").Value = 20148.428
$ws.Range("I124").Value = 9600
$ws.Range("K124").Value = 28800
$ws.Range("M124").Value = -23890

$ws.Range("H129").Value = 4306187
$ws.Range("I129").Value = 9000739
$ws.Range("K129").Value = 27002217
$ws.Range("M129").Value = -26997217

$ws.Range("H131").Value = 16753.295
$ws.Range("J131").Value = 3311.4614
$ws.Range("L131").Value = 9934.3842
$ws.Range("N131").Value = -20014.3842

$ws.Range("H135").Value = 3539.8333
$ws.Range("J135").Value = 5958.077
$ws.Range("L135").Value = 53622.693
$ws.Range("N135").Value = -58692.693

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 25577.715
$ws.Range("I126").Value = 31189.818
$ws.Range("K126").Value = 93569.454
$ws.Range("M126").Value = -91099.454

$ws.Range("H132").Value = 11250.019
$ws.Range("I132").Value = 13347.068
$ws.Range("J132").Value = 2023
$ws.Range("K132").Value = 40041.204
$ws.Range("L132").Value = 6069
$ws.Range("M132").Value = -37511.204
$ws.Range("N132").Value = -11129

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 19334
$ws.Range("I7").Value = 24501
$ws.Range("J7").Value = 9000
$ws.Range("K7").Value = 24501
$ws.Range("L7").Value = 9000
$ws.Range("M7").Value = -24389
$ws.Range("N7").Value = -9224

$ws.Range("H16").Value = 3380.4243
$ws.Range("I16").Value = 3544.8215
$ws.Range("J16").Value = 2459.8
$ws.Range("K16").Value = 3544.8215
$ws.Range("L16").Value = 2459.8
$ws.Range("M16").Value = -3374.8215
$ws.Range("N16").Value = -2799.8

$ws.Range("H22").Value = 1620.5454
$ws.Range("J22").Value = 1740
$ws.Range("L22").Value = 1740
$ws.Range("N22").Value = -2330

$ws.Range("H27").Value = 1620.5454
$ws.Range("J27").Value = 1740
$ws.Range("L27").Value = 1740
$ws.Range("N27").Value = -1954

$ws.Range("H40").Value = 7713.857
$ws.Range("I40").Value = 5999.6665
$ws.Range("K40").Value = 5999.6665
$ws.Range("M40").Value = -5863.6665

$ws.Range("H46").Value = 3277.353
$ws.Range("I46").Value = 1544
$ws.Range("K46").Value = 1544
$ws.Range("M46").Value = -1356

$ws.Range("H63").Value = 39042.5
$ws.Range("I63").Value = 38000
$ws.Range("J63").Value = 40085
$ws.Range("K63").Value = 38000
$ws.Range("L63").Value = 40085
$ws.Range("M63").Value = -37251
$ws.Range("N63").Value = -41583

$ws.Range("H66").Value = 39042.5
$ws.Range("I66").Value = 38000
$ws.Range("J66").Value = 40085
$ws.Range("K66").Value = 114000
$ws.Range("L66").Value = 120255
$ws.Range("M66").Value = -110256
$ws.Range("N66").Value = -127743

$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

$ws.Range("H122").Value = 562025.4399999999
$ws.Range("I122").Value = 913786.4399999999
$ws.Range("K122").Value = 2741359.32
$ws.Range("M122").Value = -2738909.32

$ws.Range("H126").Value = 19334
$ws.Range("I126").Value = 24501
$ws.Range("J126").Value = 9000
$ws.Range("K126").Value = 73503
$ws.Range("L126").Value = 27000
$ws.Range("M126").Value = -71033
$ws.Range("N126").Value = -31940

$ws.Range("H132").Value = 4364.5
$ws.Range("I132").Value = 3561.9375
$ws.Range("K132").Value = 10685.8125
$ws.Range("M132").Value = -8155.8125

$ws.Range("H136").Value = 4036.5264
$ws.Range("I136").Value = 4121.5625
$ws.Range("K136").Value = 12364.6875
$ws.Range("M136").Value = -9814.6875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5626.1943
$ws.Range("I126").Value = 5716.6665
$ws.Range("J126").Value = 4631
$ws.Range("K126").Value = 17149.9995
$ws.Range("L126").Value = 13893
$ws.Range("M126").Value = -14679.9995
$ws.Range("N126").Value = -18833

$ws.Range("H132").Value = 2700.7026
$ws.Range("I132").Value = 2719.0908
$ws.Range("J132").Value = 2549
$ws.Range("K132").Value = 8157.2724
$ws.Range("L132").Value = 7647
$ws.Range("M132").Value = -5627.2724
$ws.Range("N132").Value = -12707

$ws.Range("H135").Value = 54049.168
$ws.Range("J135").Value = 54049.168
$ws.Range("L135").Value = 54049.168
$ws.Range("N135").Value = -64189.168
